$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to retain their original text representation
# (some "prices" look numeric, e.g. "235.64", and Excel would otherwise
# silently coerce them into real numbers when .Value is assigned).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.215.25'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '1.855.58'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '235.64'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D7').Value = '0.4779'
$ws.Range('E7').Value = '  -2.69%  '
$ws.Range('E8').Value = '  -4.13%  '
$ws.Range('D9').Value = '0.06462'
$ws.Range('D10').Value = '1.857.39'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').Value = '0.07364'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('D13').Value = '5.081'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').Value = '87.06'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').Value = '0.6454'
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('D16').Value = '30.157.95'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '13.13'
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('D19').Value = '0.000007583'
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('D20').Value = '224.95'
$ws.Range('E20').Value = '  +16.44%  '
$ws.Range('D21').Value = '2.099.26'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = '5.279'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').Value = '6.068'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').Value = '9.194'
$ws.Range('E25').Value = '  -2.87%  '
$ws.Range('D26').Value = '163.29'
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('D28').Value = '1.922'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '4.232'
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.09160'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').Value = '3.947'
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('D33').Value = '0.04962'
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('D34').Value = '0.7325'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  +4.23%  '
$ws.Range('D36').Value = '2.687'
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').Value = '0.01837'
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('D38').Value = '2.595'
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('D39').Value = '0.8985'
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('D40').Value = '2.047'
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').Value = '5.945'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').Value = '106.01'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('D44').Value = '0.4229'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').Value = '7.358'
$ws.Range('E45').Value = '  -2.82%  '
$ws.Range('D46').Value = '0.1311'
$ws.Range('E46').Value = '  -3.77%  '
$ws.Range('D47').Value = '64.26'
$ws.Range('E47').Value = '  -6.04%  '
$ws.Range('D48').Value = '1.512'
$ws.Range('E48').Value = '  +8.03%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '33.89'
$ws.Range('E49').Value = '  -2.87%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.672'
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('D51').Value = '0.05656'
$ws.Range('E51').Value = '  -3.26%  '

# Restore default (no explicit) style on the price column so the saved
# worksheet XML matches the original cell styling (no "s" attribute).
$ws.Range("D2:D51").Style = "Normal"
